# Apply the edit described by the commit: rename sheets, tweak a handful of
# model inputs, and update the selection/view state to match.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (formula references like Sheet1!P2 auto-update to Main!P2) ---
$wsMain = $wb.Worksheets.Item("Sheet1")
$wsMain.Name = "Main"

$wsModel = $wb.Worksheets.Item("Sheet2")
$wsModel.Name = "Model"

# --- Main sheet: update Price input (P2) ---
$wsMain.Range("P2").Value = 68.5
$wsMain.Range("P2").Select()

# --- Model sheet: update hardcoded/overridden inputs ---

# J4: was a formula (=I4*0.96); now a typed-in hardcoded number.
$wsModel.Range("J4").Value = 7791

# K19: was a formula (=K17-K18); now a typed-in hardcoded number.
$wsModel.Range("K19").Value = 1240

# Q31: plain input flips sign.
$wsModel.Range("Q31").Value = -3328

# U31: was a formula (=T31*1.05); now a typed-in hardcoded number.
$wsModel.Range("U31").Value = 1900

# --- Restore view/selection state on the Model sheet ---
$wsModel.Activate()
$window = $excel.ActiveWindow
$window.Panes.Item(1).ScrollColumn = 17  # topLeftCell Q11 -> column Q
$wsModel.Range("Y27").Select()

$excel.ActiveWindow.WindowState = $excel.ActiveWindow.WindowState
